$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 596.6667
$ws.Range("I2").Value = 596.6667
$ws.Range("K2").Value = 596.6667
$ws.Range("M2").Value = -483.6667

$ws.Range("H45").Value = 8887.5
$ws.Range("J45").Value = 8887.5
$ws.Range("L45").Value = 26662.5
$ws.Range("N45").Value = -27046.5

$ws.Range("H94").Value = 1137.7778
$ws.Range("I94").Value = 1236.4286
$ws.Range("K94").Value = 1236.4286
$ws.Range("M94").Value = -785.4286

$ws.Range("H100").Value = 2377.625
$ws.Range("I100").Value = 2224.5
$ws.Range("K100").Value = 2224.5
$ws.Range("M100").Value = -1683.5

$ws.Range("H101").Value = 317
$ws.Range("I101").Value = 317
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 951
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 671
$ws.Range("N101").ClearContents()

$ws.Range("H111").Value = 905.3333
$ws.Range("J111").Value = 721.5
$ws.Range("L111").Value = 2164.5
$ws.Range("N111").Value = -8298.5

$ws.Range("H132").Value = 8169.231
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 245.66667
$ws.Range("I5").Value = 288.8
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 288.8
$ws.Range("L5").Value = 30
$ws.Range("M5").Value = -176.8
$ws.Range("N5").Value = -254

$ws.Range("H32").Value = 3945.2
$ws.Range("I32").Value = 3626.5789
$ws.Range("K32").Value = 3626.5789
$ws.Range("M32").Value = -3339.5789

$ws.Range("H37").Value = 4534
$ws.Range("I37").Value = 4534
$ws.Range("K37").Value = 4534
$ws.Range("M37").Value = -4261

$ws.Range("H122").Value = 1515.2693
$ws.Range("I122").Value = 1515.2693
$ws.Range("K122").Value = 4545.8079
$ws.Range("M122").Value = -2095.8079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 245.66667
$ws.Range("I4").Value = 288.8
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 288.8
$ws.Range("L4").Value = 30
$ws.Range("M4").Value = -173.8
$ws.Range("N4").Value = -260

$ws.Range("H20").Value = 2624.6155
$ws.Range("I20").Value = 1718.2
$ws.Range("K20").Value = 1718.2
$ws.Range("M20").Value = -1471.2

$ws.Range("H76").Value = 18333
$ws.Range("J76").Value = 18333
$ws.Range("L76").Value = 18333
$ws.Range("N76").Value = -18963

$ws.Range("H79").Value = 18333
$ws.Range("J79").Value = 18333
$ws.Range("L79").Value = 18333
$ws.Range("N79").Value = -20517

$ws.Range("H86").Value = 3560.8667
$ws.Range("I86").Value = 3954.923
$ws.Range("J86").Value = 999.5
$ws.Range("K86").Value = 3954.923
$ws.Range("L86").Value = 999.5
$ws.Range("M86").Value = -2831.923
$ws.Range("N86").Value = -3245.5

$ws.Range("H89").Value = 3560.8667
$ws.Range("I89").Value = 3954.923
$ws.Range("J89").Value = 999.5
$ws.Range("K89").Value = 19774.615
$ws.Range("L89").Value = 4997.5
$ws.Range("M89").Value = -14158.615
$ws.Range("N89").Value = -16229.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws.Range("H31").Value = 1565.2727
$ws.Range("I31").Value = 1521.8
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 1521.8
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -1226.8
$ws.Range("N31").Value = -2590

$ws.Range("H34").Value = 1565.2727
$ws.Range("I34").Value = 1521.8
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 1521.8
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = -1319.8
$ws.Range("N34").Value = -2404

$ws.Range("H45").Value = 5000
$ws.Range("J45").Value = 5000
$ws.Range("L45").Value = 5000
$ws.Range("N45").Value = -6186

$ws.Range("H99").Value = 3072.9333
$ws.Range("I99").Value = 2945.111
$ws.Range("J99").Value = 3264.6667
$ws.Range("K99").Value = 2945.111
$ws.Range("L99").Value = 3264.6667
$ws.Range("M99").Value = -1447.111
$ws.Range("N99").Value = -6260.6667

$ws.Range("H122").Value = 4339.6875
$ws.Range("I122").Value = 4860.1665
$ws.Range("K122").Value = 14580.4995
$ws.Range("M122").Value = -12130.4995

$ws.Range("H126").Value = 3072.9333
$ws.Range("I126").Value = 2945.111
$ws.Range("J126").Value = 3264.6667
$ws.Range("K126").Value = 8835.332999999999
$ws.Range("L126").Value = 9794.000100000001
$ws.Range("M126").Value = -6365.332999999999
$ws.Range("N126").Value = -14734.0001

$ws.Range("H132").Value = 2459.7778
$ws.Range("I132").Value = 2505.4285
$ws.Range("K132").Value = 7516.2855
$ws.Range("M132").Value = -4986.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4326857
$ws.Range("I4").Value = 47999.832
$ws.Range("K4").Value = 143999.496
$ws.Range("M4").Value = -143887.496

$ws.Range("H12").Value = 292.3846
$ws.Range("I12").Value = 281.5
$ws.Range("J12").Value = 309.8
$ws.Range("K12").Value = 844.5
$ws.Range("L12").Value = 929.4000000000001
$ws.Range("M12").Value = -671.5
$ws.Range("N12").Value = -1275.4

$ws.Range("H15").Value = 233.11111
$ws.Range("J15").Value = 95
$ws.Range("L15").Value = 285
$ws.Range("N15").Value = -565

$ws.Range("H47").Value = 1863.5
$ws.Range("I47").Value = 1066
$ws.Range("J47").Value = 4256
$ws.Range("K47").Value = 3198
$ws.Range("L47").Value = 12768
$ws.Range("M47").Value = -2767
$ws.Range("N47").Value = -13630

$ws.Range("H57").Value = 17199.8
$ws.Range("J57").Value = 17199.8
$ws.Range("L57").Value = 51599.39999999999
$ws.Range("N57").Value = -52717.39999999999

$ws.Range("H92").Value = 212.4
$ws.Range("I92").Value = 212.4
$ws.Range("K92").Value = 637.2
$ws.Range("M92").Value = 610.8

$ws.Range("H97").Value = 986.5
$ws.Range("I97").Value = 800.6
$ws.Range("J97").Value = 1172.4
$ws.Range("K97").Value = 2401.8
$ws.Range("L97").Value = 3517.2
$ws.Range("M97").Value = -1905.8
$ws.Range("N97").Value = -4509.200000000001

$ws.Range("H113").Value = 2230.5715
$ws.Range("I113").Value = 2899.5
$ws.Range("J113").Value = 2119.0833
$ws.Range("K113").Value = 8698.5
$ws.Range("L113").Value = 6357.249899999999
$ws.Range("M113").Value = -6528.5
$ws.Range("N113").Value = -10697.2499

$ws.Range("H141").Value = 6416.3335
$ws.Range("J141").Value = 5000
$ws.Range("L141").Value = 15000
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 10000
$ws.Range("J52").Value = 10000
$ws.Range("L52").Value = 10000
$ws.Range("N52").Value = -10518

$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 10000
$ws.Range("N58").Value = -10554
$ws.Range("M58").ClearContents()

$ws.Range("H122").Value = 2082.1428
$ws.Range("I122").Value = 2020.3529
$ws.Range("K122").Value = 6061.0587
$ws.Range("M122").Value = -3611.0587

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H126").Value = 2349.5
$ws.Range("I126").Value = 2359.8
$ws.Range("J126").Value = 2298
$ws.Range("K126").Value = 7079.400000000001
$ws.Range("L126").Value = 6894
$ws.Range("M126").Value = -4609.400000000001
$ws.Range("N126").Value = -11834

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4541.7036
$ws.Range("I7").Value = 4261
$ws.Range("K7").Value = 4261
$ws.Range("M7").Value = -4149

$ws.Range("H22").Value = 979.9
$ws.Range("I22").Value = 359.8
$ws.Range("J22").Value = 1600
$ws.Range("K22").Value = 359.8
$ws.Range("L22").Value = 1600
$ws.Range("M22").Value = -64.80000000000001
$ws.Range("N22").Value = -2190

$ws.Range("H25").Value = 4500
$ws.Range("I25").Value = 4500
$ws.Range("K25").Value = 4500
$ws.Range("M25").Value = -4270

$ws.Range("H27").Value = 979.9
$ws.Range("I27").Value = 359.8
$ws.Range("J27").Value = 1600
$ws.Range("K27").Value = 359.8
$ws.Range("L27").Value = 1600
$ws.Range("M27").Value = -252.8
$ws.Range("N27").Value = -1814

$ws.Range("H46").Value = 1416.6666
$ws.Range("I46").Value = 1425.25
$ws.Range("J46").Value = 1409.8
$ws.Range("K46").Value = 1425.25
$ws.Range("L46").Value = 1409.8
$ws.Range("M46").Value = -1237.25
$ws.Range("N46").Value = -1785.8

$ws.Range("H55").Value = 688.7619
$ws.Range("J55").Value = 899.7692
$ws.Range("L55").Value = 899.7692
$ws.Range("N55").Value = -1245.7692

$ws.Range("H122").Value = 3579.0476
$ws.Range("I122").Value = 3311.2666
$ws.Range("K122").Value = 9933.799800000001
$ws.Range("M122").Value = -7483.799800000001

$ws.Range("H126").Value = 4541.7036
$ws.Range("I126").Value = 4261
$ws.Range("K126").Value = 12783
$ws.Range("M126").Value = -10313

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 5000
$ws.Range("J13").Value = 5000
$ws.Range("L13").Value = 5000
$ws.Range("N13").Value = -5280

$ws.Range("H126").Value = 4810.375
$ws.Range("I126").Value = 4781.6665
$ws.Range("J126").Value = 4896.5
$ws.Range("K126").Value = 14344.9995
$ws.Range("L126").Value = 14689.5
$ws.Range("M126").Value = -11874.9995
$ws.Range("N126").Value = -19629.5
